$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the user role for row 3 (E3) to the new value "Supervisor"
# This adds a new shared string entry since "Supervisor" doesn't exist yet
$ws.Range("E3").Value = "Supervisor"

# Update the selected/active cell to H9 (reflects cursor position after bulk import)
$ws.Range("H9").Select()
